$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.043.25'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '2.291.91'
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.96%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.609'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0912'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.28'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.55%  '
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('E14').Value = '  -3.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.17%  '
$ws.Range('D16').Value = '2.633.88'
$ws.Range('E16').Value = '  -2.58%  '
$ws.Range('D17').Value = '2.289.71'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').Value = '41.960.27'
$ws.Range('E18').Value = '  -1.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000105'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '256.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.17%  '
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.55'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0891'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.77%  '
$ws.Range('E33').Value = '  -6.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.27%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.119'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.51%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.130'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0353'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.81'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.62%  '
$ws.Range('E40').Value = '  -4.48%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.48'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.80%  '
$ws.Range('B43').Value = 'BitcoinSV'
$ws.Range('C43').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '97.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.228'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.66%  '
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '112.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.85%  '
$ws.Range('E51').Value = '  -0.85%  '
